$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data to fill in for rows 34-46 (columns A=Remessa, B=Material, C=Quantidade)
$data = @(
    @("80266695", "10000-LDG-I", 1),
    @("80266697", "21487-MET-I", 2160),
    @("80266697", "40343-TDK-N", 1800),
    @("80266697", "15079-TDK-N", 3500),
    @("80266697", "15403-TDK-N", 200),
    @("80266697", "40353-TDK-N", 200),
    @("80266697", "60240-STM-I", 300),
    @("80266697", "40046-TDK-I", 800),
    @("80266700", "10361-ARI-I", 1),
    @("80266701", "10255-ARI-I", 1),
    @("80266702", "10399-ARI-I", 1),
    @("80266703", "10650-ARI-I", 1),
    @("80266704", "20041-CTY-I", 1)
)

$startRow = 34

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $aVal = $data[$i][0]
    $bVal = $data[$i][1]
    $cVal = $data[$i][2]

    $aCell = $ws.Cells.Item($row, 1)
    $bCell = $ws.Cells.Item($row, 2)
    $cCell = $ws.Cells.Item($row, 3)

    if ($aVal -match '^[0-9]+$') {
        # Purely-numeric codes (e.g. "80266695") need to be forced to text so
        # they keep being stored as shared strings, same as the rest of the
        # sheet. Using a TEXT() formula and pasting the computed value back
        # keeps the original cell style (no new number format is created).
        $aCell.Formula = '=TEXT(' + $aVal + ',"0")'
        $aCell.Copy() | Out-Null
        $aCell.PasteSpecial(-4163) | Out-Null # xlPasteValues
        $excel.CutCopyMode = $false
    } else {
        $aCell.Value = $aVal
    }

    if ($bVal -match '^[0-9]+$') {
        $bCell.Formula = '=TEXT(' + $bVal + ',"0")'
        $bCell.Copy() | Out-Null
        $bCell.PasteSpecial(-4163) | Out-Null # xlPasteValues
        $excel.CutCopyMode = $false
    } else {
        $bCell.Value = $bVal
    }

    $cCell.Value = $cVal
}

# Update the selection to match the diff (A2:C46 selected, active cell A2)
$ws.Range("A2:C46").Select()
